$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 10125.94
$ws.Range("B9").Value = 10185.01
$ws.Range("C9").Value = 307.21
$ws.Range("D9").Value = 305.42
$ws.Range("E9").Value = $false
$ws.Range("F9").Value = -0.58
$ws.Range("G9").Value = 42609.488993055558
$ws.Range("G8").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("H9").Value = $false
